$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, shifting existing rows 45-69 down to 46-70.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new record's data.
$ws.Range("A45").Value = 1
$ws.Range("B45").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C45").Value = "Arica y Parinacota"
$ws.Range("D45").Value = 44879
$ws.Range("D45").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E45").Value = 15
$ws.Range("F45").Value = 100112031
$ws.Range("G45").Value = "Poroto verde"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 1300
$ws.Range("K45").Value = 1300
$ws.Range("L45").Value = 1400
$ws.Range("M45").Value = 1350
$ws.Range("N45").Value = "$/kilo"
$ws.Range("O45").Value = "Región de Arica y Parinacota"
$ws.Range("P45").Value = 1350
$ws.Range("Q45").Value = 1
$ws.Range("R45").Value = "Hortaliza"
